$d = $word.ActiveDocument

# 1. Replace the {jabatan_orang_1} placeholder with the actual position
#    title "Wali Nagari" (pulled from Firestore per the commit message).
$d.Content.Find.Execute("{jabatan_orang_1}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Wali Nagari", 2)

# 2. Remove the stray leading "AN " run that sits in front of
#    "WALI NAGARI LIMO KOTO" near the signature block, turning
#    "AN WALI NAGARI LIMO KOTO" into "WALI NAGARI LIMO KOTO".
#    We find the exact character range of "AN " and delete it directly
#    (rather than doing a Find/Replace spanning both runs) so the
#    neighbouring run ("WALI NAGARI LIMO KOTO") is left completely
#    untouched, matching the original run structure.
$findRng = $d.Content
$found = $findRng.Find.Execute("AN WALI NAGARI LIMO KOTO", $true, $false, $false, $false, $false, $true)
if ($found) {
    $anRng = $d.Range($findRng.Start, $findRng.Start + 3)
    if ($anRng.Text -eq "AN ") {
        $anRng.Delete()
    }
}
